# ------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1. Insert a brand-new worksheet named "2022-Q3" right before the
#    existing "2022-Q2" tab (so the tab order becomes:
#    总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3).
#    We create it by duplicating "2022-Q2" (so it inherits identical
#    formatting/column layout) and then replace its data with the
#    2022-Q3 fund holdings.
# 2. Refresh the "总计" (totals) summary sheet so it gets a new first
#    data row for 2022-Q3 and the existing rows keep following below.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Build the "2022-Q3" sheet ---------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Fund holdings for 2022-Q3 (basic code, name, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名)
$q3rows = @(
  @(0, "506007", "广发科创板两年定开混合",            "5.01", "94.25", "7.96", "0.3988", 1),
  @(1, "010296", "万家互联互通中国优势量化策略混合A", "4.22", "86.53", "5.27", "0.2224", 9),
  @(2, "320012", "诺安主题精选混合",                  "3.79", "86.75", "4.64", "0.1759", 5),
  @(3, "310388", "申万菱信消费增长混合A",              "2.87", "91.17", "4.32", "0.1240", 5),
  @(4, "013000", "广发盛泽一年持有期混合A",            "2.29", "82.39", "3.54", "0.0811", 9),
  @(5, "002133", "广发鑫益灵活配置混合",                "1.22", "93.85", "4.95", "0.0604", 7),
  @(6, "010297", "万家互联互通中国优势量化策略混合C", "0.46", "86.53", "5.27", "0.0242", 9),
  @(7, "006881", "华宝大健康混合",                    "0.76", "91.19", "2.71", "0.0206", 10),
  @(8, "013001", "广发盛泽一年持有期混合C",            "0.27", "82.39", "3.54", "0.0096", 9),
  @(9, "015254", "申万菱信消费增长混合C",              "0.05", "91.17", "4.32", "0.0022", 5)
)

for ($i = 0; $i -lt $q3rows.Count; $i++) {
    $row = 2 + $i
    $rec = $q3rows[$i]

    if ($row -gt 2) {
        # new row: clone the formatting of row 2's index cell (style "2")
        # onto column A of this row before writing its value.
        $q3.Cells.Item(2, 1).Copy()
        $q3.Cells.Item($row, 1).PasteSpecial(-4122)
    }

    $q3.Cells.Item($row, 1).Value = $rec[0]

    $c2 = $q3.Cells.Item($row, 2)
    $c2.NumberFormat = "@"
    $c2.Value = $rec[1]

    $q3.Cells.Item($row, 3).Value = $rec[2]

    $c4 = $q3.Cells.Item($row, 4)
    $c4.NumberFormat = "@"
    $c4.Value = $rec[3]

    $c5 = $q3.Cells.Item($row, 5)
    $c5.NumberFormat = "@"
    $c5.Value = $rec[4]

    $c6 = $q3.Cells.Item($row, 6)
    $c6.NumberFormat = "@"
    $c6.Value = $rec[5]

    $c7 = $q3.Cells.Item($row, 7)
    $c7.NumberFormat = "@"
    $c7.Value = $rec[6]

    $q3.Cells.Item($row, 8).Value = $rec[7]
}

# ---- 2. Refresh the "总计" summary sheet --------------------------
$zj = $wb.Worksheets.Item("总计")

# clone the index-column formatting (style "2") onto the new row 6
$zj.Cells.Item(5, 1).Copy()
$zj.Cells.Item(6, 1).PasteSpecial(-4122)

$zjrows = @(
  @(0, "2022-Q3", 10, 1.12),
  @(1, "2022-Q2", 1, 0.3),
  @(2, "2022-Q1", 1, 0),
  @(3, "2021-Q4", 3, 0.01),
  @(4, "2021-Q3", 4, 0.02)
)

for ($i = 0; $i -lt $zjrows.Count; $i++) {
    $row = 2 + $i
    $rec = $zjrows[$i]
    $zj.Cells.Item($row, 1).Value = $rec[0]
    $zj.Cells.Item($row, 2).Value = $rec[1]
    $zj.Cells.Item($row, 3).Value = $rec[2]
    $zj.Cells.Item($row, 4).Value = $rec[3]
}
